$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows 2-13: numeroprocesso, setoratual, datacadastro(serial), localizacao_caixa, responsavel
$data = @(
    @{ Row = 2;  A = "2024/193641.3"; B = "SECRETARIA 2-B"; C = 45634; D = "SANDERLAN"; E = "Fulano2" },
    @{ Row = 3;  A = "2024/152522.1"; B = "SECRETARIA 2-B"; C = 45602; D = "MARCOS";    E = "Fulano2" },
    @{ Row = 4;  A = "2024/113415.2"; B = "SECRETARIA 2-B"; C = 45397; D = "SABRINA";   E = "Fulano1`n" },
    @{ Row = 5;  A = "2024/193641.4"; B = "SECRETARIA 2-B"; C = 45634; D = "SANDERLAN"; E = "Fulano6" },
    @{ Row = 6;  A = "2024/152522.5"; B = "SECRETARIA 2-B"; C = 45602; D = "MARCOS";    E = "Fulano3" },
    @{ Row = 7;  A = "2024/173425.7"; B = "SECRETARIA 2-B"; C = 45397; D = "SABRINA";   E = "Fulano6" },
    @{ Row = 8;  A = "2024/193821.3"; B = "SECRETARIA 2-B"; C = 45634; D = "SANDERLAN"; E = "Fulano3" },
    @{ Row = 9;  A = "2024/152522.8"; B = "SECRETARIA 2-B"; C = 45602; D = "MARCOS";    E = "Fulano5" },
    @{ Row = 10; A = "2024/113415.5"; B = "SECRETARIA 2-B"; C = 45397; D = "SABRINA";   E = "Fulano4`n" },
    @{ Row = 11; A = "2024/196641.3"; B = "SECRETARIA 2-B"; C = 45634; D = "SANDERLAN"; E = "Fulano4`n" },
    @{ Row = 12; A = "2024/122522.9"; B = "SECRETARIA 2-B"; C = 45602; D = "MARCOS";    E = "Fulano1`n" },
    @{ Row = 13; A = "2024/119605.2"; B = "SECRETARIA 2-B"; C = 45397; D = "SABRINA";   E = "Fulano5" }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
}

# New rows 12 and 13 need the same date number format (style index 2 = YYYY-MM-DD) as existing column C cells
$ws.Cells.Item(12, 3).NumberFormat = $ws.Cells.Item(11, 3).NumberFormat
$ws.Cells.Item(13, 3).NumberFormat = $ws.Cells.Item(11, 3).NumberFormat
